$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataEntry")

# --- Define new text values as here-strings (exact newlines / trailing spaces preserved) ---
$B6 = @'
FM radio stations, and possibly favorited stations. 
Database Used: Firebase
'@
$B7 = @'
Setting or selecting different FM radio stations, Favorite different FM radio stations, Displays song or station currently playing 
'@
$B8 = @'
SparkFun Electronics (For Parts), Elmwood Electronics (For Parts), Prototype Lab (For Extra Help)
'@
$B10 = @'
This project which we will be creating will be able to connect to a speaker via FM Radio. How this  will work is by taking a mobile device and connecting to the database in order for the FM Radio stations to play through the speaker bonnet.
'@
$B11 = @'
One of the products being used will be the FM evaluation board tuner chip. This device does more then tuning into FM stations, it can also detect both data service and radio broadcast data service. It can also be used to display station id and song to the user as well as have great filtering and carrying detection. This board will be able to pick up multiple radio stations and makes a great tool in order for it to be implemented with a Raspberry Pi. By using a speaker bonnet that acts as the output for this sensor, that is the primary source in which the sound will be coming from. It will amplify the audio so it can play in areas that you usually cannot hear with your mobile speakers.
'@
$B12 = @'
Bluetooth in wireless communication. (n.d.). Retrieved September 18, 2017, from http://ieeexplore.ieee.org/document/1007414/
SparkFun FM Tuner Evaluation Board - Si4703. (n.d.). Retrieved from https://www.sparkfun.com/products/12938

'@
$B13 = @'
Lumpkins, W. (n.d.). The MobiAria Wireless Bluetooth Speaker. Retrieved September 18, 2017, from http://ieeexplore.ieee.org/document/6685931/
Bodson, D. (n.d.). Digital Audio Around the World. Retrieved from http://ieeexplore.ieee.org/stamp/stamp.jsp?tp=&arnumber=5641649
Pauli, M. (2017, May 5). Miniaturized Millimeter-Wave Radar Sensor for High-Accuracy Applications. Retrieved from http://ieeexplore.ieee.org/stamp/stamp.jsp?tp=&arnumber=7885501               
'@
$B15 = @'
Our end solution and goal for this project is to use an amplified speaker and output the radio signal stations from the FM Radio sensor.
'@

# --- Rich text field B14 (bold headers + normal descriptions) ---
$run1 = @'
CanaKit Raspberry Pi 3 Complete Starter Kit - 32 GB Edition -
'@
$run2 = @'
 The platform in which we will be doing this project on.

'@
$run3 = @'
SPARKFUN FM TUNER EVALUATION BOARD - SI4703 - 
'@
$run4 = @'
Enables users to tune into FM radio stations.

'@
$run5 = @'
JUMPER WIRES - CONNECTED 6" (M/F, 20 PACK) - 
'@
$run6 = @'
Used for connection between breadboard where sensor is attached to Raspberry Pi 3.
'@
$run7 = @'

Adafruit I2S 3W Stereo Speaker Bonnet for Raspberry Pi - Mini Kit - 
'@
$run8 = @'
Speaker used for amplyifing audio from the Raspberry Pi 3. 
'@
$full14 = $run1 + $run2 + $run3 + $run4 + $run5 + $run6 + $run7 + $run8

# --- Write cells in the same order the original workbook author touched them ---
$ws.Range("B15").Value = $B15
$ws.Range("B13").Value = $B13
$ws.Range("B12").Value = $B12
$ws.Range("B11").Value = $B11
$ws.Range("B10").Value = $B10
$ws.Range("B8").Value = $B8
$ws.Range("B7").Value = $B7
$ws.Range("B6").Value = $B6
$ws.Range("B14").Value = $full14

# Apply bold formatting to the header portion of each run in B14
$pos = 1
$r = $ws.Range("B14").Characters($pos, $run1.Length)
$r.Font.Bold = $true
$pos = $pos + $run1.Length

$r = $ws.Range("B14").Characters($pos, $run2.Length)
$r.Font.Bold = $false
$pos = $pos + $run2.Length

$r = $ws.Range("B14").Characters($pos, $run3.Length)
$r.Font.Bold = $true
$pos = $pos + $run3.Length

$r = $ws.Range("B14").Characters($pos, $run4.Length)
$r.Font.Bold = $false
$pos = $pos + $run4.Length

$r = $ws.Range("B14").Characters($pos, $run5.Length)
$r.Font.Bold = $true
$pos = $pos + $run5.Length

$r = $ws.Range("B14").Characters($pos, $run6.Length)
$r.Font.Bold = $false
$pos = $pos + $run6.Length

$r = $ws.Range("B14").Characters($pos, $run7.Length)
$r.Font.Bold = $true
$pos = $pos + $run7.Length

$r = $ws.Range("B14").Characters($pos, $run8.Length)
$r.Font.Bold = $false
$pos = $pos + $run8.Length

# --- Match the final selection left behind in the source workbook ---
$ws.Activate() | Out-Null
$ws.Range("A15").Select() | Out-Null

